$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header "TD-HBN"
$ws.Range("F1").Value = "TD-HBN"

# Rename B3's text: "Validation Accuracy" -> "Validation" + newline + "Accuracy", center+wrap style
$ws.Range("B3").Value = "Validation`nAccuracy"
$ws.Range("B3").WrapText = $true

# Column widths (character widths equivalent to stored widths 12.28515625 / 60.28515625)
$ws.Columns.Item(2).ColumnWidth = 11.570870535714286
$ws.Columns.Item(6).ColumnWidth = 59.570870535714285

# Update selection
$ws.Range("F3").Select() | Out-Null
